$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count") entirely — this shifts F:K left to E:J,
# matching the new header order reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date.
$ws.Columns("E").Delete()
